# Continued Work on Card Management
# - finalize the CSV-database-style layout on Sheet1:
#   * header row: drop "defacto"/"range" strings, rename C1 header to "sub_type",
#     swap G1/H1 headers to action_type_percent_rate / pop_change_rate
#   * row 2 (Snowmagdon / weather) gets a sub_type + numeric pop_change_rate
#   * new row 3 (Ice System 1), row 4 (Ocotpi) and row 5 (Earth) are added
#   * a few column widths are widened to fit the new, longer content
#   * view settings (zoom + active cell) are updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "type"
$ws.Range("C1").Value = "sub_type"
$ws.Range("D1").Value = "title"
$ws.Range("E1").Value = "description"
$ws.Range("F1").Value = "population"
$ws.Range("G1").Value = "action_type_percent_rate"
$ws.Range("H1").Value = "pop_change_rate"
$ws.Range("I1").Value = "image_location"

# ---- Row 2 (existing "Snowmagdon" weather card) ----
$ws.Range("B2").Value = "action"
$ws.Range("C2").Value = "weather"
$ws.Range("D2").Value = "Snowmagdon"
$ws.Range("E2").Value = "Let is snow!"
$ws.Range("F2").Value = -100
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = 0.3

# ---- Row 3 (new "Ice System 1" weather card) ----
$ws.Range("B3").Value = "solar"
$ws.Range("C3").Value = "weather"
$ws.Range("D3").Value = "Ice System 1"
$ws.Range("E3").Value = "We like ice in our solar system:)"
$ws.Range("F3").Value = 300
$ws.Range("G3").Value = 0.25
$ws.Range("H3").Value = 0.3

# ---- Row 4 (new "Ocotpi" solar card) ----
$ws.Range("B4").Value = "solar"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "Ocotpi"
$ws.Range("E4").Value = "We are an allien race of octopi!"
$ws.Range("F4").Value = 500
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = 0.2

# ---- Row 5 (new "Earth" planet card) ----
$ws.Range("B5").Value = "planet"
$ws.Range("C5").Value = "null"
$ws.Range("D5").Value = "Earth"
$ws.Range("E5").Value = "We are human."
$ws.Range("F5").Value = 600
$ws.Range("G5").Value = "null"
$ws.Range("H5").Value = 0.03

# ---- Column widths: widen the columns that now hold longer text ----
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666   # title column
$ws.Columns.Item(7).ColumnWidth = 28.666666666666668   # action_type_percent_rate column (best-fit)
$ws.Columns.Item(8).ColumnWidth = 20.333333333333332   # pop_change_rate column

# ---- View: zoom out a bit and land the selection on the last entry ----
[void]$ws.Range("H5").Select()
$excel.ActiveWindow.Zoom = 140
